$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (preserve rich-text run formatting via Characters) ---
$ws.Range("A8").Characters(21,1).Text = "8"
$ws.Range("C9").Characters(27,9).Text = "2/17/2025"
$ws.Range("C9").Characters(47,9).Text = "2/23/2025"

# --- Pure numeric value updates (style/type unchanged) ---
$ws.Range("H15").Value = -66.666666666666
$ws.Range("I15").Value = 2
$ws.Range("J15").Value = 6
$ws.Range("K15").Value = -66.666666666666
$ws.Range("L15").Value = -33.333333333333
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = 100
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 60
$ws.Range("I16").Value = 20
$ws.Range("K16").Value = 53.846153846153
$ws.Range("L16").Value = 150
$ws.Range("M16").Value = -23.076923076923
$ws.Range("N16").Value = -84.496124031007
$ws.Range("C17").Value = 1
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 12
$ws.Range("H17").Value = -14.285714285714
$ws.Range("I17").Value = 20
$ws.Range("J17").Value = 26
$ws.Range("K17").Value = -23.076923076923
$ws.Range("L17").Value = -31.03448275862
$ws.Range("M17").Value = 25
$ws.Range("N17").Value = -42.857142857142
$ws.Range("C18").Value = 3
$ws.Range("I18").Value = 20
$ws.Range("K18").Value = 11.111111111111
$ws.Range("L18").Value = -4.761904761904
$ws.Range("M18").Value = -67.213114754098
$ws.Range("N18").Value = -93.150684931506
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 38
$ws.Range("H19").Value = -21.052631578947
$ws.Range("I19").Value = 55
$ws.Range("J19").Value = 84
$ws.Range("K19").Value = -34.523809523809
$ws.Range("L19").Value = -40.217391304347
$ws.Range("M19").Value = -5.172413793103
$ws.Range("N19").Value = -47.115384615384
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -66.666666666666
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = -28.571428571428
$ws.Range("I20").Value = 20
$ws.Range("J20").Value = 28
$ws.Range("K20").Value = -28.571428571428
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -4.761904761904
$ws.Range("N20").Value = -92.88256227758
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = -6.25
$ws.Range("F21").Value = 73
$ws.Range("G21").Value = 82
$ws.Range("H21").Value = -10.975609756097
$ws.Range("I21").Value = 137
$ws.Range("J21").Value = 177
$ws.Range("K21").Value = -22.598870056497
$ws.Range("L21").Value = -21.264367816092
$ws.Range("M21").Value = -25.136612021857
$ws.Range("N21").Value = -83.767772511848
$ws.Range("H22").Value = -100
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = 5.882352941176
$ws.Range("F24").Value = 115
$ws.Range("G24").Value = 84
$ws.Range("H24").Value = 36.904761904761
$ws.Range("I24").Value = 176
$ws.Range("J24").Value = 145
$ws.Range("K24").Value = 21.379310344827
$ws.Range("L24").Value = 2.923976608187
$ws.Range("M24").Value = 22.222222222222
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 17
$ws.Range("H25").Value = -10.526315789473
$ws.Range("I25").Value = 26
$ws.Range("J25").Value = 34
$ws.Range("K25").Value = -23.529411764705
$ws.Range("L25").Value = -49.019607843137
$ws.Range("C26").Value = 8
$ws.Range("E26").Value = 14.285714285714
$ws.Range("G26").Value = 28
$ws.Range("H26").Value = 3.571428571428
$ws.Range("I26").Value = 52
$ws.Range("J26").Value = 56
$ws.Range("K26").Value = -7.142857142857
$ws.Range("L26").Value = 33.333333333333
$ws.Range("M26").Value = -10.344827586206
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 3
$ws.Range("J27").Value = 8
$ws.Range("K27").Value = -62.5
$ws.Range("L27").Value = -25
$ws.Range("F28").Value = 5
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 8
$ws.Range("K28").Value = -38.461538461538
$ws.Range("L28").Value = 60

# --- Text -> Number conversions (set target number format, then numeric value) ---
$ws.Range("C15").NumberFormat = '#,##0'
$ws.Range("C15").Value = 1
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("D15").Value = 1
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E15").Value = 0
$ws.Range("F15").NumberFormat = '#,##0'
$ws.Range("F15").Value = 1
$ws.Range("C27").NumberFormat = '#,##0'
$ws.Range("C27").Value = 2
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("D27").Value = 1
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E27").Value = 100
$ws.Range("F27").NumberFormat = '#,##0'
$ws.Range("F27").Value = 2

# --- Number -> Text conversions (force text storage, then copy the canonical text style) ---
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "0"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "***.*"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "***.*"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "***.*"
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = "0"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"

# Re-apply the canonical text-cell style (General numFmt, right/center aligned, Andale WT 9.1)
$ws.Range("C14").Copy() | Out-Null
$ws.Range("G14,H14,D16,E16,D18,E18,F22,D28,E28").PasteSpecial(-4122)
